$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain decimal numbers (e.g. "4.10", "0.0499")
# must be forced to Text format first, otherwise Excel auto-converts them to
# numeric values and the exact display string (trailing zeros, leading zeros,
# etc.) would be lost.
$textCells = @(
    "D5", "D9", "D14", "D15", "D16", "D18", "D19", "D25", "D30", "D31", "D32", "D34", "D37",
    "D38", "D42", "D43", "D44", "D47", "D50"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values scraped for this run of the GitHub Actions job.
$ws.Range("D2").Value = '27.037.09'
$ws.Range("E2").Value = '  +0.71%  '
$ws.Range("D3").Value = '1.679.50'
$ws.Range("E3").Value = '  +1.00%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '216.10'
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  +1.73%  '
$ws.Range("D9").Value = '21.34'
$ws.Range("E9").Value = '  +5.68%  '
$ws.Range("E10").Value = '  +0.71%  '
$ws.Range("E11").Value = '  -0.81%  '
$ws.Range("D12").Value = '1.916.93'
$ws.Range("E12").Value = '  +1.04%  '
$ws.Range("D13").Value = '1.726.53'
$ws.Range("E13").Value = '  +3.84%  '
$ws.Range("D14").Value = '4.10'
$ws.Range("E14").Value = '  +0.56%  '
$ws.Range("D15").Value = '0.533'
$ws.Range("E15").Value = '  +1.89%  '
$ws.Range("D16").Value = '66.50'
$ws.Range("E16").Value = '  +0.81%  '
$ws.Range("D17").Value = '27.042.96'
$ws.Range("E17").Value = '  +0.75%  '
$ws.Range("D18").Value = '8.17'
$ws.Range("E18").Value = '  +3.34%  '
$ws.Range("D19").Value = '235.61'
$ws.Range("E19").Value = '  +1.35%  '
$ws.Range("E20").Value = '  +1.04%  '
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("E22").Value = '  +1.10%  '
$ws.Range("E23").Value = '  +1.21%  '
$ws.Range("E24").Value = '  -4.02%  '
$ws.Range("D25").Value = '146.59'
$ws.Range("E25").Value = '  +0.51%  '
$ws.Range("E26").Value = '  +1.74%  '
$ws.Range("E27").Value = '  +3.59%  '
$ws.Range("E28").Value = '  -2.44%  '
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("D30").Value = '0.0499'
$ws.Range("E30").Value = '  +0.72%  '
$ws.Range("D31").Value = '1.18'
$ws.Range("E31").Value = '  +0.35%  '
$ws.Range("D32").Value = '3.37'
$ws.Range("E32").Value = '  +0.65%  '
$ws.Range("D33").Value = '1.533.11'
$ws.Range("E33").Value = '  +5.04%  '
$ws.Range("D34").Value = '3.18'
$ws.Range("E34").Value = '  +0.97%  '
$ws.Range("E35").Value = '  +5.02%  '
$ws.Range("E36").Value = '  -0.79%  '
$ws.Range("D37").Value = '0.589'
$ws.Range("E37").Value = '  +2.63%  '
$ws.Range("D38").Value = '0.921'
$ws.Range("E38").Value = '  +2.46%  '
$ws.Range("E39").Value = '  +3.38%  '
$ws.Range("E40").Value = '  +6.48%  '
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '68.05'
$ws.Range("E42").Value = '  +3.59%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '5.61'
$ws.Range("E43").Value = '  -2.71%  '
$ws.Range("D44").Value = '2.25'
$ws.Range("E44").Value = '  -0.11%  '
$ws.Range("D45").Value = '1.822.75'
$ws.Range("E45").Value = '  +0.72%  '
$ws.Range("E46").Value = '  +0.46%  '
$ws.Range("D47").Value = '90.35'
$ws.Range("E47").Value = '  +0.09%  '
$ws.Range("E48").Value = '  +0.18%  '
$ws.Range("E49").Value = '  +2.22%  '
$ws.Range("D50").Value = '8.01'
$ws.Range("E50").Value = '  +6.35%  '
$ws.Range("E51").Value = '  -0.20%  '
